# Scheduled-runner data refresh: recompute profit-tracking figures on the
# Seraph_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Only the
# computed H:N value columns change; row identity (A:G) is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2500
$ws.Range("H98").Value = 620.7778
$ws.Range("I98").Value = 689
$ws.Range("K98").Value = 689
$ws.Range("M98").Value = 809
$ws.Range("H113").Value = 3290.75
$ws.Range("I113").Value = 2582
$ws.Range("K113").Value = 2582
$ws.Range("M113").Value = 672
$ws.Range("H122").Value = 620.7778
$ws.Range("I122").Value = 689
$ws.Range("K122").Value = 2067
$ws.Range("M122").Value = 383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1627.25
$ws.Range("I45").Value = 1627.25
$ws.Range("K45").Value = 1627.25
$ws.Range("M45").Value = -1250.25
$ws.Range("H61").Value = 3244.6667
$ws.Range("I61").Value = 3244.6667
$ws.Range("K61").Value = 3244.6667
$ws.Range("M61").Value = -3032.6667
$ws.Range("H136").Value = 3244.6667
$ws.Range("I136").Value = 3244.6667
$ws.Range("K136").Value = 9734.000100000001
$ws.Range("M136").Value = -7184.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3814.111
$ws.Range("I105").Value = 3137.1
$ws.Range("J105").Value = 5748.4287
$ws.Range("K105").Value = 3137.1
$ws.Range("L105").Value = 5748.4287
$ws.Range("M105").Value = -1390.1
$ws.Range("N105").Value = -9242.4287
$ws.Range("H134").Value = 2369.16
$ws.Range("I134").Value = 2162.1738
$ws.Range("J134").Value = 4749.5
$ws.Range("K134").Value = 6486.5214
$ws.Range("L134").Value = 14248.5
$ws.Range("M134").Value = -3951.5214
$ws.Range("N134").Value = -19318.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 941.3077
$ws.Range("I16").Value = 860.7273
$ws.Range("K16").Value = 860.7273
$ws.Range("M16").Value = -573.7273
$ws.Range("H86").Value = 12206.5
$ws.Range("I86").Value = 9447.333000000001
$ws.Range("K86").Value = 9447.333000000001
$ws.Range("M86").Value = -8324.333000000001
$ws.Range("H89").Value = 12206.5
$ws.Range("I89").Value = 9447.333000000001
$ws.Range("K89").Value = 47236.665
$ws.Range("M89").Value = -41620.665
$ws.Range("H99").Value = 10299.728
$ws.Range("J99").Value = 12676.053
$ws.Range("L99").Value = 12676.053
$ws.Range("N99").Value = -15672.053
$ws.Range("H105").Value = 2107.318
$ws.Range("I105").Value = 525.9286
$ws.Range("K105").Value = 525.9286
$ws.Range("M105").Value = 1221.0714
$ws.Range("H113").Value = 941.3077
$ws.Range("I113").Value = 860.7273
$ws.Range("K113").Value = 860.7273
$ws.Range("M113").Value = 1309.2727
$ws.Range("H126").Value = 10299.728
$ws.Range("J126").Value = 12676.053
$ws.Range("L126").Value = 38028.159
$ws.Range("N126").Value = -42968.159

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 46432524
$ws.Range("J4").Value = 5498.5
$ws.Range("L4").Value = 16495.5
$ws.Range("N4").Value = -16719.5
$ws.Range("H34").Value = 894.1111
$ws.Range("I34").Value = 726.3333
$ws.Range("J34").Value = 978
$ws.Range("K34").Value = 2178.9999
$ws.Range("L34").Value = 2934
$ws.Range("M34").Value = -2094.9999
$ws.Range("N34").Value = -3102
$ws.Range("H62").Value = 1349.75
$ws.Range("J62").Value = 1300
$ws.Range("L62").Value = 3900
$ws.Range("N62").Value = -5272
$ws.Range("H65").Value = 1349.75
$ws.Range("J65").Value = 1300
$ws.Range("L65").Value = 11700
$ws.Range("N65").Value = -18564
$ws.Range("H131").Value = 2105.6191
$ws.Range("I131").Value = 1329.909
$ws.Range("J131").Value = 2958.9
$ws.Range("K131").Value = 3989.727
$ws.Range("L131").Value = 8876.700000000001
$ws.Range("M131").Value = 1050.273
$ws.Range("N131").Value = -18956.7
$ws.Range("H132").Value = 1549.3846
$ws.Range("I132").Value = 1543
$ws.Range("J132").Value = 1563.75
$ws.Range("K132").Value = 13887
$ws.Range("L132").Value = 14073.75
$ws.Range("M132").Value = -11357
$ws.Range("N132").Value = -19133.75
$ws.Range("H140").Value = 2741.077
$ws.Range("I140").Value = 2741.077
$ws.Range("K140").Value = 8223.231
$ws.Range("M140").Value = -3043.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 26666
$ws.Range("J52").Value = 26666
$ws.Range("L52").Value = 26666
$ws.Range("N52").Value = -27184
$ws.Range("H70").Value = 7347.231
$ws.Range("I70").Value = 5703
$ws.Range("K70").Value = 5703
$ws.Range("M70").Value = -5433
$ws.Range("H73").Value = 7347.231
$ws.Range("I73").Value = 5703
$ws.Range("K73").Value = 5703
$ws.Range("M73").Value = -4767

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1986.7142
$ws.Range("I7").Value = 2650
$ws.Range("J7").Value = 1721.4
$ws.Range("K7").Value = 2650
$ws.Range("L7").Value = 1721.4
$ws.Range("M7").Value = -2538
$ws.Range("N7").Value = -1945.4
$ws.Range("H16").Value = 1360.5333
$ws.Range("I16").Value = 1360.5333
$ws.Range("K16").Value = 1360.5333
$ws.Range("M16").Value = -1190.5333
$ws.Range("H55").Value = 564.75
$ws.Range("I55").Value = 420
$ws.Range("K55").Value = 420
$ws.Range("M55").Value = -247
$ws.Range("H126").Value = 1986.7142
$ws.Range("I126").Value = 2650
$ws.Range("J126").Value = 1721.4
$ws.Range("K126").Value = 7950
$ws.Range("L126").Value = 5164.200000000001
$ws.Range("M126").Value = -5480
$ws.Range("N126").Value = -10104.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1142.25
$ws.Range("I81").Value = 1142.25
$ws.Range("K81").Value = 2284.5
$ws.Range("M81").Value = -1223.5
$ws.Range("H84").Value = 1142.25
$ws.Range("I84").Value = 1142.25
$ws.Range("K84").Value = 11422.5
$ws.Range("M84").Value = -6118.5
$ws.Range("H96").Value = 1265.3334
$ws.Range("I96").Value = 1265.3334
$ws.Range("K96").Value = 1265.3334
$ws.Range("M96").Value = 107.6666
$ws.Range("H107").Value = 657.94116
$ws.Range("I107").Value = 458.5
$ws.Range("K107").Value = 1375.5
$ws.Range("M107").Value = 544.5
$ws.Range("H122").Value = 2838.4443
$ws.Range("I122").Value = 1069.5
$ws.Range("K122").Value = 3208.5
$ws.Range("M122").Value = -758.5
$ws.Range("H126").Value = 4474.25
$ws.Range("I126").Value = 948.5
$ws.Range("K126").Value = 2845.5
$ws.Range("M126").Value = -375.5
$ws.Range("H132").Value = 965.7879
$ws.Range("I132").Value = 965.7879
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2897.3637
$ws.Range("M132").Value = -367.3636999999999

# WVR!N132 is dropped entirely by the refresh (row 132 no longer reports a
# shipping-cost-adjusted total), so clear it rather than leave a stale value.
$wb.Worksheets.Item("WVR").Range("N132").ClearContents()
